$wb = $excel.ActiveWorkbook
$win = $wb.Windows.Item(1)
$win | Get-Member
